$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{Row=2; D='286.50'; E='2.52%'; G='22'},
    @{Row=3; E='4.30%'; G='22'},
    @{Row=4; E='4.56%'; G='22'},
    @{Row=5; D='0.06647'; E='4.11%'; G='22'},
    @{Row=6; D='7.363'; E='4.69%'; G='22'},
    @{Row=7; D='3.393'; E='1.88%'; G='22'},
    @{Row=8; D='1.378'; E='5.65%'; G='22'},
    @{Row=9; D='0.9404'; E='4.11%'; G='22'},
    @{Row=10; D='0.1566'; E='2.21%'; G='22'},
    @{Row=11; D='0.06619'; E='7.99%'; G='22'},
    @{Row=12; D='0.07614'; E='1.07%'; G='22'},
    @{Row=13; D='0.02945'; E='-0.36%'; G='22'},
    @{Row=14; D='0.09005'; E='-0.11%'; G='22'},
    @{Row=15; D='0.001609'; E='1.85%'; G='22'},
    @{Row=16; D='0.04481'; E='1.28%'; G='22'},
    @{Row=17; D='0.0006455'; E='0.25%'; G='22'},
    @{Row=18; D='0.006316'; E='4.59%'; G='22'},
    @{Row=19; D='3.456'; E='-0.90%'; G='22'},
    @{Row=20; D='2.251'; E='0.96%'; G='22'},
    @{Row=21; D='0.3215'; E='2.24%'; G='22'},
    @{Row=22; E='-4.08%'; G='22'},
    @{Row=23; D='4.085'; E='4.43%'; G='22'},
    @{Row=24; D='0.1554'; E='3.22%'; G='22'},
    @{Row=25; D='0.001183'; E='0.43%'; G='22'},
    @{Row=26; D='0.004477'; E='4.62%'; G='22'},
    @{Row=27; D='0.0001251'; E='6.00%'; G='22'},
    @{Row=28; D='0.0001620'; E='-2.42%'; G='22'},
    @{Row=29; G='22'},
    @{Row=30; G='22'},
    @{Row=31; G='22'},
    @{Row=32; G='22'},
    @{Row=33; G='22'},
    @{Row=34; G='22'},
    @{Row=35; G='22'},
    @{Row=36; G='22'},
    @{Row=37; G='22'},
    @{Row=38; G='22'},
    @{Row=39; G='22'},
    @{Row=40; D='0.04212'; E='3.42%'; G='22'},
    @{Row=41; D='0.006765'; E='1.37%'; G='22'},
    @{Row=42; D='0.1255'; E='-10.92%'; G='22'},
    @{Row=43; D='0.002021'; E='-2.82%'; G='22'},
    @{Row=44; D='0.01226'; E='11.46%'; G='22'},
    @{Row=45; D='0.00005708'; E='2.95%'; G='22'},
    @{Row=46; E='25.93%'; G='22'},
    @{Row=47; D='0.01308'; E='-29.35%'; G='22'},
    @{Row=48; G='22'},
    @{Row=49; G='22'},
    @{Row=50; G='22'},
    @{Row=51; G='22'},
)

$colMap = @{ D=4; E=5; G=7 }

foreach ($item in $changes) {
    $r = $item.Row
    foreach ($col in @('D','E','G')) {
        if ($item.ContainsKey($col)) {
            $cell = $ws.Cells.Item($r, $colMap[$col])
            $cell.NumberFormat = "@"
            $cell.Value = $item[$col]
        }
    }
}
